$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 173, pushing the existing rows 173..210 down to 174..211.
$ws.Rows.Item(173).Insert()

# Populate the newly inserted row 173 with the new weekly record.
$ws.Range("A173").Value = 5
$ws.Range("B173").Value = "Macroferia Regional de Talca"
$ws.Range("C173").Value = "Maule"
$ws.Range("D173").Value = "3/17/2022"
$ws.Range("E173").Value = 7
$ws.Range("F173").Value = 100112045
$ws.Range("G173").Value = "Zapallo"
$ws.Range("H173").Value = "Camote"
$ws.Range("I173").Value = "1a (cosecha)"
$ws.Range("J173").Value = 900
$ws.Range("K173").Value = 300
$ws.Range("L173").Value = 300
$ws.Range("M173").Value = 300
$ws.Range("N173").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O173").Value = "Región del Maule"
$ws.Range("P173").Value = 300
$ws.Range("Q173").Value = 1
$ws.Range("R173").Value = "Hortaliza"
